$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated save_data now records K (strike count-derived value) instead of Strike#
# in column G. Update the per-row values accordingly.
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 1
